$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.037.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.223.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.95%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.512"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.60%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.42%  "

$ws.Range("E11").Value = "  -2.99%  "

$ws.Range("E12").Value = "  +2.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.46"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.567.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.221.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.727"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "39.970.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0884"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.12%  "

$ws.Range("E28").Value = "  -1.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "156.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.68%  "

$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0716"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.33%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.26%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.67%  "

$ws.Range("E37").Value = "  +0.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0974"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.117.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("E43").Value = "  -4.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.03%  "

$ws.Range("E45").Value = "  -1.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.433.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.17%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "88.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.41%  "
